$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the row-2 formatting pattern onto the new row 9 first ---
# (row 2 is the other "normal"/6-line-script row: D col = hyperlink style,
#  E col = wrap style, which is exactly the visual pattern row 9 needs)
$ws.Range("A2:E2").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)

# --- Now write the new row's content. Write order matters: it controls the
# order new entries land in sharedStrings.xml (target has E's script text at
# index 27, A's label at 28, D's url at 29) ---
$ws.Range("E9").Value = "url_list = soup.select(""#table tr td.tal"")[:5];`nurl_list = [""https://www.hanaw.com/main/research/trends/RC_060600_P1.cmd"" + i.select_one(""a"")[""href""]for i in url_list];`ntitle_list = soup.select(""#table tr td.tal a.title"")[:5];`ntitle_list = [i.text.strip() for i in title_list];`nbody_list = soup.select(""#table tr td.tdlast"")[:5];`nbody_list = [i.text.strip() for i in body_list];"
$ws.Range("A9").Value = "#하나증권"
$ws.Range("D9").Value = "https://www.hanaw.com/main/research/trends/RC_060600_P1.cmd"
$ws.Range("B9").Value = 1.1
$ws.Range("C9").Value = "normal"

# Row height to match the other 6-line wrapped rows (ht=99, same as rows 2/3/6/7)
$ws.Range("A9:E9").RowHeight = 99

# Hyperlink for D9, pointing at the hanaw research URL (text-to-display already
# matches thanks to the Value set above, like the other D-column hyperlinks)
$ws.Hyperlinks.Add($ws.Range("D9"), "https://www.hanaw.com/main/research/trends/RC_060600_P1.cmd")

# Hyperlinks.Add stamps its own cell style onto D9; restore the normal
# hyperlink-cell style (same as D2/D3/.../D6) that the rest of the column uses
$ws.Range("D2").Copy()
$ws.Range("D9").PasteSpecial(-4122)

# Move the active selection to E4 (matches author's last-saved cursor position)
$ws.Range("E4").Select()
